$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("F1").Value = "MotivoAnulacion"
$ws.Range("F2").Value = "Anulación por Desistimiento"

$ws.Range("F3").Select()
